$wb = $excel.ActiveWorkbook
$dst = $wb.Worksheets.Item("Rough")
for ($i = 0; $i -le 9; $i++) {
  $cell = $dst.Cells.Item(100 + $i, 26)
  $cell.Value = "x"
  $cell.Interior.Pattern = 1
  $cell.Interior.ThemeColor = $i
}
